# Insert a new daily price record for "Ciboulette" at row 333 of the
# "Vega Central Mapocho de Santiago" sheet. Excel's EntireRow.Insert()
# shifts the existing row 333 (and everything below it) down by one row,
# growing the used range from A1:R397 to A1:R398, and the new blank row
# inherits the date-formatted style from its neighbours (column D).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(333).EntireRow.Insert()

# Fill the newly inserted row with the new record's data. Every column
# except D (Fecha) and J (Volumen) reuses the values that used to sit in
# (old) row 333, since the new record shares the same market/category/
# quality/price-range/origin metadata as the record that followed it.
$ws.Cells.Item(333, 1).Value = 9
$ws.Cells.Item(333, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(333, 3).Value = "Metropolitana"
$ws.Cells.Item(333, 4).Value = 44694
$ws.Cells.Item(333, 5).Value = 13
$ws.Cells.Item(333, 6).Value = 100112039
$ws.Cells.Item(333, 7).Value = "Ciboulette"
$ws.Cells.Item(333, 8).Value = "Sin especificar"
$ws.Cells.Item(333, 9).Value = "Primera"
$ws.Cells.Item(333, 10).Value = 250
$ws.Cells.Item(333, 11).Value = 800
$ws.Cells.Item(333, 12).Value = 1000
$ws.Cells.Item(333, 13).Value = 900
$ws.Cells.Item(333, 14).Value = "`$/docena de atados"
$ws.Cells.Item(333, 15).Value = "Región Metropolitana"
$ws.Cells.Item(333, 16).Value = 300
$ws.Cells.Item(333, 17).Value = 3
$ws.Cells.Item(333, 18).Value = "Hortaliza"
